$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make sure the right sheet is active (it already is the selected tab)
$ws.Activate()

# --- Update data values (new regression data) ---
$ws.Range("E2").Value = 3199801311
$ws.Range("N2").Value = 30990137

# --- Update the view: scroll position and active selection ---
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("N3").Select()
